$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column - copy header style (bold, bordered, centered) from G1
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "Save" values per row (2-41), derived from whether G > 9 (era/sum threshold)
$saveValues = @(
    0,0,0,1,0,0,1,0,0,0,
    0,0,0,0,0,0,1,0,0,0,
    0,1,1,0,0,0,0,0,1,0,
    0,0,0,0,0,1,0,0,0,0
)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
